$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B22").Value = 6292
$ws.Range("C22").Value = 989
$ws.Range("D22").Value = 5772949
$ws.Range("E22").Value = 917.5061983471074
$ws.Range("F22").Value = 8.314684110862448
$ws.Range("G22").Value = 3.451882845188292
$ws.Range("H22").Value = 25.54365034490558
